$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had 16 data rows (2-17); the refreshed exposure-site
# data now only needs 15 data rows (2-16), so drop the now-unused last row.
$ws.Rows.Item(17).Delete()

# Rewrite every data cell with the refreshed Case Alerts exposure-site table.
$ws.Cells.Item(2,1).Value2 = 'Broadmeadows'
$ws.Cells.Item(2,2).Value2 = 'Craigieburn Line train'
$ws.Cells.Item(2,3).Value2 = '1:25pm - 1:59pm  9/02/2021'
$ws.Cells.Item(2,4).Value2 = 'Case caught train from Broadmeadows Railway Station to Glenroy Railway Station'
$ws.Cells.Item(2,5).Value2 = 'old'

$ws.Cells.Item(3,1).Value2 = 'Broadmeadows'
$ws.Cells.Item(3,2).Value2 = 'Craigieburn Line train'
$ws.Cells.Item(3,3).Value2 = '1:25pm - 1:59pm  9/2/2021'
$ws.Cells.Item(3,4).Value2 = 'Case caught train from Broadmeadows Railway Station to Glenroy Railway Station'
$ws.Cells.Item(3,5).Value2 = 'new'

$ws.Cells.Item(4,1).Value2 = 'Broadmeadows'
$ws.Cells.Item(4,2).Value2 = 'Sacca''s Fruit World  Broadmeadows Central  Broadmeadows VIC 3047'
$ws.Cells.Item(4,3).Value2 = '12:30pm - 1:00pm  9/2/2021'
$ws.Cells.Item(4,4).Value2 = 'Case visited venue'
$ws.Cells.Item(4,5).Value2 = 'new'

$ws.Cells.Item(5,1).Value2 = 'Broadmeadows'
$ws.Cells.Item(5,2).Value2 = 'Woolworths  Broadmeadows Central  Pascoe Vale Road  Broadmeadows VIC 3047'
$ws.Cells.Item(5,3).Value2 = '12:15pm - 12:30 pm 9/2/2021'
$ws.Cells.Item(5,4).Value2 = 'Case attended venue'
$ws.Cells.Item(5,5).Value2 = 'old'

$ws.Cells.Item(6,1).Value2 = 'Broadmeadows'
$ws.Cells.Item(6,2).Value2 = 'Woolworths  Broadmeadows Central  Pascoe Vale Road  Broadmeadows VIC 3047'
$ws.Cells.Item(6,3).Value2 = '12:15pm - 12:30pm 9/2/2021'
$ws.Cells.Item(6,4).Value2 = 'Case attended venue'
$ws.Cells.Item(6,5).Value2 = 'new'

$ws.Cells.Item(7,1).Value2 = 'Coburg'
$ws.Cells.Item(7,2).Value2 = 'Function venue  426 Sydney Rd  Coburg VIC 3058'
$ws.Cells.Item(7,3).Value2 = '7:14pm  11:30pm  6/02/2021'
$ws.Cells.Item(7,4).Value2 = 'Case attended venue'
$ws.Cells.Item(7,5).Value2 = 'old'

$ws.Cells.Item(8,1).Value2 = 'Coburg'
$ws.Cells.Item(8,2).Value2 = 'Function venue  426 Sydney Rd  Coburg VIC 3058'
$ws.Cells.Item(8,3).Value2 = '7:14pm  11:30pm  6/2/2021'
$ws.Cells.Item(8,4).Value2 = 'Case attended venue'
$ws.Cells.Item(8,5).Value2 = 'new'

$ws.Cells.Item(9,1).Value2 = 'Glenroy'
$ws.Cells.Item(9,2).Value2 = '513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham'
$ws.Cells.Item(9,3).Value2 = '1:35pm  2:17pm  9/02/2021'
$ws.Cells.Item(9,4).Value2 = 'Case caught bus from Glenroy Railway Station towards Eltham'
$ws.Cells.Item(9,5).Value2 = 'old'

$ws.Cells.Item(10,1).Value2 = 'Glenroy'
$ws.Cells.Item(10,2).Value2 = '513 Eltham to Glenroy bus route  Glenroy Railway Station towards Eltham'
$ws.Cells.Item(10,3).Value2 = '1:35pm  2:17pm  9/2/2021'
$ws.Cells.Item(10,4).Value2 = 'Case caught bus from Glenroy Railway Station towards Eltham'
$ws.Cells.Item(10,5).Value2 = 'new'

$ws.Cells.Item(11,1).Value2 = 'Heatherton'
$ws.Cells.Item(11,2).Value2 = 'Melbourne Golf Academy  385 Centre Dandenong Rd  Heatherton VIC 3202'
$ws.Cells.Item(11,3).Value2 = '5:19pm - 6:35pm  1/2/2021'
$ws.Cells.Item(11,4).Value2 = 'Case attended venue'
$ws.Cells.Item(11,5).Value2 = 'new'

$ws.Cells.Item(12,1).Value2 = 'Heatherton'
$ws.Cells.Item(12,2).Value2 = 'Melbourne Golf Academy  385 Centre Dandenong Rd  Heatherton, VIC 3202'
$ws.Cells.Item(12,3).Value2 = '5:19pm - 6:35pm  1/2/2021'
$ws.Cells.Item(12,4).Value2 = 'Case attended venue'
$ws.Cells.Item(12,5).Value2 = 'old'

$ws.Cells.Item(13,1).Value2 = 'Hoppers Crossing'
$ws.Cells.Item(13,2).Value2 = 'Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing VIC 3029'
$ws.Cells.Item(13,3).Value2 = '6:45am - 7:30am  8/02/21'
$ws.Cells.Item(13,4).Value2 = 'Case attended venue'
$ws.Cells.Item(13,5).Value2 = 'old'

$ws.Cells.Item(14,1).Value2 = 'Hoppers Crossing'
$ws.Cells.Item(14,2).Value2 = 'Coates Hire Werribee  148A Geelong Rd  Hoppers Crossing VIC 3029'
$ws.Cells.Item(14,3).Value2 = '6:45am - 7:30am  8/2/21'
$ws.Cells.Item(14,4).Value2 = 'Case attended venue'
$ws.Cells.Item(14,5).Value2 = 'new'

$ws.Cells.Item(15,1).Value2 = 'Melbourne'
$ws.Cells.Item(15,2).Value2 = 'Queen Victoria Market  Queen Street  Melbourne VIC 3000'
$ws.Cells.Item(15,3).Value2 = '8:25am - 10:10am  11/2/2021'
$ws.Cells.Item(15,4).Value2 = 'Case attended Section 2 - Fruit and Vegetables, and used Section 2 female toilets. See a map of the Queen Victoria Market (PDF)'
$ws.Cells.Item(15,5).Value2 = 'old'

$ws.Cells.Item(16,1).Value2 = 'Melbourne'
$ws.Cells.Item(16,2).Value2 = 'Queen Victoria Market  Queen Street  Melbourne VIC 3000'
$ws.Cells.Item(16,3).Value2 = '8:25am - 10:10am  11/2/2021'
$ws.Cells.Item(16,4).Value2 = 'Case attended sheds A and B (also known as section 2) - Fruit and Vegetables, and used female toilets adjacent to shed A.  See a map of the Queen Victoria Market (PDF)'
$ws.Cells.Item(16,5).Value2 = 'new'

# Columns B (Site) and D (Notes) were resized to fit the new, narrower/
# wider text after the content refresh.
$ws.Columns.Item(2).ColumnWidth = 65
$ws.Columns.Item(4).ColumnWidth = 134.5

# Leave the selection where the author last left it.
$ws.Range("B16").Select() | Out-Null
